$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (pushes current rows 54-57 down to 55-58)
$ws.Rows.Item(54).Insert()

# Copy the date cell style (s="2" -> custom date format) from D53 into the new D54
$ws.Range("D53").Copy() | Out-Null
$ws.Range("D54").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Fill in the new row 54 values (same as the row below it, except date/volume)
$ws.Cells.Item(54, 1).Value = 4
$ws.Cells.Item(54, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(54, 3).Value = "Los Lagos"
$ws.Cells.Item(54, 4).Value = 44516
$ws.Cells.Item(54, 5).Value = 10
$ws.Cells.Item(54, 6).Value = 100112026
$ws.Cells.Item(54, 7).Value = "Haba"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 140
$ws.Cells.Item(54, 11).Value = 10000
$ws.Cells.Item(54, 12).Value = 10000
$ws.Cells.Item(54, 13).Value = 10000
$ws.Cells.Item(54, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(54, 15).Value = "Región del Maule"
$ws.Cells.Item(54, 16).Value = 400
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"
